$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Beta) values C2:N2
$ws.Range("C2").Value = 19.12075701903682
$ws.Range("D2").Value = 0.006944391349456487
$ws.Range("E2").Value = 0.01982943797740053
$ws.Range("F2").Value = 34.41294661881737
$ws.Range("G2").Value = 0.01653725991617913
$ws.Range("H2").Value = 105.6611796974577
$ws.Range("I2").Value = 0.0004032156549604761
$ws.Range("J2").Value = [double]"2.651477957110355e-07"
$ws.Range("K2").Value = 0.001206095338689192
$ws.Range("L2").Value = 0.006376174261821647
$ws.Range("M2").Value = [double]"2.988608620412108e-06"
$ws.Range("N2").Value = 0.01953834335799941

# Update row 3 (Gamma) values C3:N3
$ws.Range("C3").Value = 0.04981522627320694
$ws.Range("D3").Value = 0.04815098319456564
$ws.Range("E3").Value = 0.0499839736740351
$ws.Range("F3").Value = 0.04777808792960521
$ws.Range("G3").Value = 0.04747010325951689
$ws.Range("H3").Value = 0.04810047702392242
$ws.Range("I3").Value = 0.04618026579441559
$ws.Range("J3").Value = 0.04588368811761902
$ws.Range("K3").Value = 0.04648968793614636
$ws.Range("L3").Value = 0.04781083512558815
$ws.Range("M3").Value = 0.04750290028369024
$ws.Range("N3").Value = 0.04813319895520735

# Add new row 4 (Beta + Gamma)
# Copy the formatting from A2 (bold, centered, bordered) onto A4, same as
# the other index cells in column A.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A4").Value = 2

$ws.Range("B4").Value = "Beta + Gamma"

$ws.Range("C4").Value = 19.17057224531003
$ws.Range("D4").Value = 0.05509537454402212
$ws.Range("E4").Value = 0.06981341165143562
$ws.Range("F4").Value = 34.46072470674698
$ws.Range("G4").Value = 0.06400736317569602
$ws.Range("H4").Value = 105.7092801744816
$ws.Range("I4").Value = 0.04658348144937609
$ws.Range("J4").Value = 0.04588395326541474
$ws.Range("K4").Value = 0.04769578327483555
$ws.Range("L4").Value = 0.05418700938740979
$ws.Range("M4").Value = 0.04750588889231065
$ws.Range("N4").Value = 0.06767154231320675
